$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 49.301159418377814
$ws.Range("C2").Value = -4.977969127720586
$ws.Range("D2").Value = 1.1676935802359607
$ws.Range("E2").Value = 8.8685582203942523

$ws.Range("B3").Value = 46.018332723408086
$ws.Range("C3").Value = 9.1606185307708188
$ws.Range("D3").Value = -13.455662020161444
$ws.Range("E3").Value = 41.709503591796292

$ws.Range("B1:E3").Select() | Out-Null
